$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D11: "Nugu Auto" how-to article text -- replace the HTML <a href="..">Nugu Auto</a>
# link markup with a bare URL (matches the sharedStrings.xml change in the diff).
$ws.Range("D11").Value = "고객이 Nugu Auto를 이용해 음성 인식 시스템에 대해 문의하는 경우 아래를 참고합니다.`n<br>`n<br>https://www.volvocars.com/kr/support/car/xc40-recharge-pure-electric/article/a419a85488f7c458c0a8015153fd7d99"

# D12: "TMAP Auto" how-to article text -- replace the HTML <a href="..">TMAP help</a>
# link markup with a bare URL (matches the sharedStrings.xml change in the diff).
$ws.Range("D12").Value = "TMAP Auto에 대해 문의하는 경우 아래 페이지를 참고합니다`n<br>`n<br>https://tmaphelp.zendesk.com/hc/ko/sections/25814297789083-%EB%B3%BC%EB%B3%B4"

# Scroll the sheet view down one row (topLeftCell moves from A7 to A8) while
# keeping the existing selection on D11.
$ws.Application.ActiveWindow.ScrollRow = 8

Write-Host "Updated D11 and D12 link text; adjusted sheet scroll position"
